# Fruta / hortaliza, semanal
# Insert one new daily price record for "Granada" at row 83 (Vega Modelo de
# Temuco). This shifts the existing rows 83-156 down to 84-157 and the new
# row 83 duplicates the data that is now on row 84, except for an updated
# date (D) and volume (M).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 83, pushing everything else down.
$ws.Rows.Item(83).Insert()

# Populate the newly inserted row 83 with the new record.
$ws.Cells.Item(83, 1).Value = 10
$ws.Cells.Item(83, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(83, 3).Value = "La Araucanía"
$ws.Cells.Item(83, 4).Value = 44790
$ws.Cells.Item(83, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(83, 5).Value = 9
$ws.Cells.Item(83, 6).Value = "Fruta"
$ws.Cells.Item(83, 7).Value = 100104
$ws.Cells.Item(83, 8).Value = "Frutos de pepita"
$ws.Cells.Item(83, 9).Value = 100104001
$ws.Cells.Item(83, 10).Value = "Granada"
$ws.Cells.Item(83, 11).Value = "Wonderfull"
$ws.Cells.Item(83, 12).Value = "Primera"
$ws.Cells.Item(83, 13).Value = 150
$ws.Cells.Item(83, 14).Value = 14000
$ws.Cells.Item(83, 15).Value = 14000
$ws.Cells.Item(83, 16).Value = 14000
$ws.Cells.Item(83, 17).Value = "`$/bandeja 10 kilos granel"
$ws.Cells.Item(83, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(83, 19).Value = 1400
$ws.Cells.Item(83, 20).Value = 10
